$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 24.38000000000037
$ws.Range("H2").Value = 0.000007254619304841725
$ws.Range("I2").Value = 0.000007254619304841725
$ws.Range("L2").Value = 62.48732254290727
$ws.Range("M2").Value = "[36.01063153784034, 88.9640135479742]"
$ws.Range("N2").Value = 0.00002080584869568192
$ws.Range("O2").Value = 0.00002080584869568192
$ws.Range("P2").Value = 1.654131867655886
$ws.Range("Q2").Value = "[1.1383949355350387, 2.169868799776734]"
$ws.Range("R2").Value = 0.00000006468368662915225
$ws.Range("S2").Value = 0.00000006468368662915225
$ws.Range("T2").Value = 81.1400125460621
$ws.Range("U2").Value = "[65.07057309192803, 97.20945200019617]"
$ws.Range("V2").Value = 0.0000000000003059774655866931
$ws.Range("W2").Value = 0.0000000000003059774655866931
$ws.Range("X2").Value = 17.96164164164192
$ws.Range("Y2").Value = 15.96048048048072
$ws.Range("Z2").Value = 19.96280280280311

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("F3").Value = 24.38000000000037
$ws.Range("H3").Value = 0.05031659429675106
$ws.Range("I3").Value = 0.05031659429675106
$ws.Range("L3").Value = 28.59408185510627
$ws.Range("M3").Value = "[-1.5246685964268565, 58.71283230663939]"
$ws.Range("N3").Value = 0.06223643482290409
$ws.Range("O3").Value = 0.06223643482290409
$ws.Range("P3").Value = 1.817658211986886
$ws.Range("Q3").Value = "[0.3836579616996545, 3.2516584622741185]"
$ws.Range("R3").Value = 0.01414240296605551
$ws.Range("S3").Value = 0.01414240296605551
$ws.Range("T3").Value = 60.02228215241456
$ws.Range("U3").Value = "[43.81433231625237, 76.23023198857675]"
$ws.Range("V3").Value = 0.000000002138252686023634
$ws.Range("W3").Value = 0.000000002138252686023634
$ws.Range("X3").Value = 17.32712712712739
$ws.Range("Y3").Value = 11.76292292292311
$ws.Range("Z3").Value = 22.89133133133168

# Row 4
$ws.Range("B4").Value = 0
$ws.Range("F4").Value = 24.38000000000037
$ws.Range("H4").Value = 0.08414785568095096
$ws.Range("I4").Value = 0.08414785568095096
$ws.Range("L4").Value = 29.26253160477975
$ws.Range("M4").Value = "[-3.895574737082775, 62.42063794664228]"
$ws.Range("N4").Value = 0.08224850435241726
$ws.Range("O4").Value = 0.08224850435241726
$ws.Range("P4").Value = 1.150973885098963
$ws.Range("Q4").Value = "[-0.6981317007977337, 3.000079470995659]"
$ws.Range("R4").Value = 0.216433801759232
$ws.Range("S4").Value = 0.216433801759232
$ws.Range("T4").Value = 79.1714577717544
$ws.Range("U4").Value = "[60.87955450474203, 97.46336103876678]"
$ws.Range("V4").Value = 0.00000000003188893593630837
$ws.Range("W4").Value = 0.00000000003188893593630837
$ws.Range("X4").Value = 19.9139939939943
$ws.Range("Y4").Value = 12.73909909909929
$ws.Range("Z4").Value = 27.08888888888931

# Row 5
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = 24.38000000000037
$ws.Range("H5").Value = 0.01549585911680429
$ws.Range("I5").Value = 0.01549585911680429
$ws.Range("L5").Value = 34.97073242771287
$ws.Range("M5").Value = "[3.5146739491241874, 66.42679090630156]"
$ws.Range("N5").Value = 0.03013529584651664
$ws.Range("O5").Value = 0.03013529584651664
$ws.Range("P5").Value = 0.6855527512338089
$ws.Range("Q5").Value = "[-0.13836844520315417, 1.509473947670772]"
$ws.Range("R5").Value = 0.1007018788451084
$ws.Range("S5").Value = 0.1007018788451084
$ws.Range("T5").Value = 70.52818163241425
$ws.Range("U5").Value = "[54.13937875123393, 86.91698451359457]"
$ws.Range("V5").Value = 0.00000000003757150146554977
$ws.Range("W5").Value = 0.00000000003757150146554977
$ws.Range("X5").Value = 21.71991991992025
$ws.Range("Y5").Value = 18.52294294294322
$ws.Range("Z5").Value = 24.91689689689728

# Row 6
$ws.Range("F6").Value = 24.38000000000037
$ws.Range("H6").Value = 0.01571506084413621
$ws.Range("I6").Value = 0.01571506084413621
$ws.Range("L6").Value = 35.70484197053132
$ws.Range("M6").Value = "[5.948188565067852, 65.46149537599479]"
$ws.Range("N6").Value = 0.01977863585759776
$ws.Range("O6").Value = 0.01977863585759776
$ws.Range("P6").Value = 1.327079178993887
$ws.Range("Q6").Value = "[0.15723686954903915, 2.4969214884387343]"
$ws.Range("R6").Value = 0.02709002042181075
$ws.Range("S6").Value = 0.02709002042181075
$ws.Range("T6").Value = 59.92761222028641
$ws.Range("U6").Value = "[43.00233849374986, 76.85288594682295]"
$ws.Range("V6").Value = 0.000000006509114802710769
$ws.Range("W6").Value = 0.000000006509114802710769
$ws.Range("X6").Value = 19.23067067067096
$ws.Range("Y6").Value = 14.69145145145168
$ws.Range("Z6").Value = 23.76988988989025

# Row 7
$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 24.38000000000037
$ws.Range("H7").Value = 0.02291512544609031
$ws.Range("I7").Value = 0.02291512544609031
$ws.Range("L7").Value = 34.67251519983571
$ws.Range("M7").Value = "[2.761987833188158, 66.58304256648327]"
$ws.Range("N7").Value = 0.03386955803572866
$ws.Range("O7").Value = 0.03386955803572866
$ws.Range("P7").Value = 1.226447582482502
$ws.Range("Q7").Value = "[-0.05660527303765406, 2.5095004380026573]"
$ws.Range("R7").Value = 0.06053288942992796
$ws.Range("S7").Value = 0.06053288942992796
$ws.Range("T7").Value = 74.86916815124047
$ws.Range("U7").Value = "[57.566973822769214, 92.17136247971173]"
$ws.Range("V7").Value = 0.00000000003211653165635653
$ws.Range("W7").Value = 0.00000000003211653165635653
$ws.Range("X7").Value = 19.62114114114144
$ws.Range("Y7").Value = 14.64264264264287
$ws.Range("Z7").Value = 24.59963963964002

# Row 8
$ws.Range("F8").Value = 24.38000000000037
$ws.Range("H8").Value = 0.0003982490472733691
$ws.Range("I8").Value = 0.0003982490472733691
$ws.Range("L8").Value = 45.00585557383052
$ws.Range("M8").Value = "[16.320177242010274, 73.69153390565076]"
$ws.Range("N8").Value = 0.002820065462665555
$ws.Range("O8").Value = 0.002820065462665555
$ws.Range("P8").Value = 1.025184389459732
$ws.Range("Q8").Value = "[0.3836579616996545, 1.6667108172198102]"
$ws.Range("R8").Value = 0.002391073876872651
$ws.Range("S8").Value = 0.002391073876872651
$ws.Range("T8").Value = 69.62117497891998
$ws.Range("U8").Value = "[54.784558948615484, 84.45779100922448]"
$ws.Range("V8").Value = 0.000000000002955635736157092
$ws.Range("W8").Value = 0.000000000002955635736157092
$ws.Range("X8").Value = 20.40208208208239
$ws.Range("Y8").Value = 17.91283283283311
$ws.Range("Z8").Value = 22.89133133133168

# Row 9
$ws.Range("F9").Value = 24.38000000000037
$ws.Range("H9").Value = 0.0009933875576764173
$ws.Range("I9").Value = 0.0009933875576764173
$ws.Range("L9").Value = 37.2654842156353
$ws.Range("M9").Value = "[13.149246059640461, 61.38172237163013]"
$ws.Range("N9").Value = 0.003221736538184583
$ws.Range("O9").Value = 0.003221736538184583
$ws.Range("P9").Value = 1.188710733790733
$ws.Range("Q9").Value = "[0.4968685077749626, 1.8805529598065025]"
$ws.Range("R9").Value = 0.001191599870456805
$ws.Range("S9").Value = 0.001191599870456805
$ws.Range("T9").Value = 48.46507639952312
$ws.Range("U9").Value = "[35.21570591391441, 61.71444688513183]"
$ws.Range("V9").Value = 0.000000002915240271406105
$ws.Range("W9").Value = 0.000000002915240271406105
$ws.Range("X9").Value = 19.76756756756787
$ws.Range("Y9").Value = 17.08308308308334
$ws.Range("Z9").Value = 22.45205205205239

# Row 10
$ws.Range("B10").Value = 1
$ws.Range("F10").Value = 24.38000000000037
$ws.Range("H10").Value = 0.00002470621480465951
$ws.Range("I10").Value = 0.00002470621480465951
$ws.Range("L10").Value = 59.42408679915083
$ws.Range("M10").Value = "[27.75886656312872, 91.08930703517294]"
$ws.Range("N10").Value = 0.0004594648103493082
$ws.Range("O10").Value = 0.0004594648103493082
$ws.Range("P10").Value = 0.761026448617347
$ws.Range("Q10").Value = "[0.2830263651882694, 1.2390265320464247]"
$ws.Range("R10").Value = 0.002473180316380885
$ws.Range("S10").Value = 0.002473180316380885
$ws.Range("T10").Value = 69.1046154404614
$ws.Range("U10").Value = "[52.89950202604841, 85.3097288548744]"
$ws.Range("V10").Value = 0.00000000004867772851468999
$ws.Range("W10").Value = 0.00000000004867772851468999
$ws.Range("X10").Value = 21.42706706706739
$ws.Range("Y10").Value = 19.57233233233263
$ws.Range("Z10").Value = 23.28180180180216

# Row 11
$ws.Range("F11").Value = 25.40000000000053
$ws.Range("H11").Value = 0.0005602004479130507
$ws.Range("I11").Value = 0.0005602004479130507
$ws.Range("L11").Value = 50.32152062967008
$ws.Range("M11").Value = "[17.487428582199342, 83.15561267714081]"
$ws.Range("N11").Value = 0.003457639622547637
$ws.Range("O11").Value = 0.003457639622547637
$ws.Range("P11").Value = 0.798763297309117
$ws.Range("Q11").Value = "[0.22013161736865516, 1.3773949772495788]"
$ws.Range("R11").Value = 0.00790019176030854
$ws.Range("S11").Value = 0.00790019176030854
$ws.Range("T11").Value = 58.48806527405323
$ws.Range("U11").Value = "[41.679086527948954, 75.29704402015751]"
$ws.Range("V11").Value = 0.000000009907907161377238
$ws.Range("W11").Value = 0.000000009907907161377238
$ws.Range("X11").Value = 22.17097097097143
$ws.Range("Y11").Value = 19.83183183183225
$ws.Range("Z11").Value = 24.51011011011062

# Row 12
$ws.Range("B12").Value = 0
$ws.Range("F12").Value = 25.40000000000053
$ws.Range("H12").Value = 0.001824923236732889
$ws.Range("I12").Value = 0.001824923236732889
$ws.Range("L12").Value = 53.19806663187404
$ws.Range("M12").Value = "[20.942898228439645, 85.45323503530844]"
$ws.Range("N12").Value = 0.001781779422746155
$ws.Range("O12").Value = 0.001781779422746155
$ws.Range("P12").Value = 0.1698158191129622
$ws.Range("Q12").Value = "[-0.6100790538502698, 0.9497106920761942]"
$ws.Range("R12").Value = 0.663082521949752
$ws.Range("S12").Value = 0.663082521949752
$ws.Range("T12").Value = 71.14593939143631
$ws.Range("U12").Value = "[51.88283456760693, 90.40904421526571]"
$ws.Range("V12").Value = 0.000000002287597222760951
$ws.Range("W12").Value = 0.000000002287597222760951
$ws.Range("X12").Value = 24.71351351351403
$ws.Range("Y12").Value = 21.56076076076121
$ws.Range("Z12").Value = 27.86626626626685

# Row 13
$ws.Range("F13").Value = 25.40000000000053
$ws.Range("H13").Value = 0.000002096627631731351
$ws.Range("I13").Value = 0.000002096627631731351
$ws.Range("L13").Value = 63.13497151216747
$ws.Range("M13").Value = "[38.081245712227044, 88.18869731210789]"
$ws.Range("N13").Value = 0.000007151166895535965
$ws.Range("O13").Value = 0.000007151166895535965
$ws.Range("P13").Value = 0.2956053147521933
$ws.Range("Q13").Value = "[-0.16981581911296129, 0.7610264486173479]"
$ws.Range("R13").Value = 0.2073739217523789
$ws.Range("S13").Value = 0.2073739217523789
$ws.Range("T13").Value = 64.11872562228882
$ws.Range("U13").Value = "[49.459623074511654, 78.77782817006599]"
$ws.Range("V13").Value = 0.00000000002357514183870535
$ws.Range("W13").Value = 0.00000000002357514183870535
$ws.Range("X13").Value = 24.20500500500551
$ws.Range("Y13").Value = 22.32352352352399
$ws.Range("Z13").Value = 26.08648648648703

# Row 14
$ws.Range("F14").Value = 25.40000000000053
$ws.Range("H14").Value = 0.008150523710598212
$ws.Range("I14").Value = 0.008150523710598212
$ws.Range("L14").Value = 37.45984201700403
$ws.Range("M14").Value = "[10.788911537016475, 64.13077249699158]"
$ws.Range("N14").Value = 0.006952959087939803
$ws.Range("O14").Value = 0.006952959087939803
$ws.Range("P14").Value = 0.05660527303765406
$ws.Range("Q14").Value = "[-0.8931054190385392, 1.0063159651138474]"
$ws.Range("R14").Value = 0.904981165556695
$ws.Range("S14").Value = 0.904981165556695
$ws.Range("T14").Value = 55.70450837327778
$ws.Range("U14").Value = "[39.85978858246037, 71.54922816409518]"
$ws.Range("V14").Value = 0.000000007731923101772509
$ws.Range("W14").Value = 0.000000007731923101772509
$ws.Range("X14").Value = 25.1711711711717
$ws.Range("Y14").Value = 21.33193193193238
$ws.Range("Z14").Value = 29.01041041041102

# Row 15
$ws.Range("F15").Value = 25.40000000000053
$ws.Range("H15").Value = 0.00005954874530744725
$ws.Range("I15").Value = 0.00005954874530744725
$ws.Range("L15").Value = 61.17966698853333
$ws.Range("M15").Value = "[31.483760148402922, 90.87557382866375]"
$ws.Range("N15").Value = 0.0001460452009212077
$ws.Range("O15").Value = 0.0001460452009212077
$ws.Range("P15").Value = -0.2515789912784623
$ws.Range("Q15").Value = "[-0.855368570346771, 0.3522105877898465]"
$ws.Range("R15").Value = 0.4057878611303565
$ws.Range("S15").Value = 0.4057878611303565
$ws.Range("T15").Value = 68.61224533655879
$ws.Range("U15").Value = "[51.227260104161104, 85.99723056895647]"
$ws.Range("V15").Value = 0.0000000004090787708577182
$ws.Range("W15").Value = 0.0000000004090787708577182
$ws.Range("X15").Value = 1.017017017017043
$ws.Range("Y15").Value = -1.42382382382385
$ws.Range("Z15").Value = 3.457857857857936

# Row 16
$ws.Range("F16").Value = 25.40000000000053
$ws.Range("H16").Value = 0.00155736335829848
$ws.Range("I16").Value = 0.00155736335829848
$ws.Range("L16").Value = 47.95053665784654
$ws.Range("M16").Value = "[13.802103329122332, 82.09896998657075]"
$ws.Range("N16").Value = 0.006965557145631962
$ws.Range("O16").Value = 0.006965557145631962
$ws.Range("P16").Value = -1.006315965113847
$ws.Range("Q16").Value = "[-1.6981581911296173, -0.3144737390980774]"
$ws.Range("R16").Value = 0.005312837190214914
$ws.Range("S16").Value = 0.005312837190214914
$ws.Range("T16").Value = 69.47525671583251
$ws.Range("U16").Value = "[51.38133340011386, 87.56918003155117]"
$ws.Range("V16").Value = 0.0000000008442140320141789
$ws.Range("W16").Value = 0.0000000008442140320141789
$ws.Range("X16").Value = 4.06806806806815
$ws.Range("Y16").Value = 1.271271271271296
$ws.Range("Z16").Value = 6.864864864865005

# Row 17
$ws.Range("F17").Value = 25.40000000000053
$ws.Range("H17").Value = 0.00001875437336751329
$ws.Range("I17").Value = 0.00001875437336751329
$ws.Range("L17").Value = 59.0858619020409
$ws.Range("M17").Value = "[33.65256195442504, 84.51916184965675]"
$ws.Range("N17").Value = 0.00002655801932949231
$ws.Range("O17").Value = 0.00002655801932949231
$ws.Range("P17").Value = -1.509473947670771
$ws.Range("Q17").Value = "[-2.07552667804731, -0.9434212172942322]"
$ws.Range("R17").Value = 0.000002650989207486631
$ws.Range("S17").Value = 0.000002650989207486631
$ws.Range("T17").Value = 59.84901033463142
$ws.Range("U17").Value = "[43.47558049915287, 76.22244017010996]"
$ws.Range("V17").Value = 0.000000002968848056283946
$ws.Range("W17").Value = 0.000000002968848056283946
$ws.Range("X17").Value = 6.102102102102229
$ws.Range("Y17").Value = 3.813813813813895
$ws.Range("Z17").Value = 8.390390390390564

